$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header D1: "Benefit" -> "Status"
$ws.Range("D1").Value = "Status"

# Column D (rows 2-15): convert raw amount to a 0/1 win-status flag,
# formatted as an integer number format.
$ws.Range("D2:D15").NumberFormat = "0"

$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0

# Page setup: letter/A4-style paper (9 = A4), portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to F4.
[void]$ws.Range("F4").Select()
